$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (losing trailing zeros / precision).
$textCells = @("D5","D6","D9","D10","D12","D16","D20","D21","D22","D23","D24","D25","D27","D28","D30","D32","D33","D34","D35","D36","D39","D41","D42","D45","D47","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.853.55"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "2.296.52"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "116.27"
$ws.Range("E5").Value = "  +18.53%  "
$ws.Range("D6").Value = "269.61"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "0.618"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("D10").Value = "48.68"
$ws.Range("E10").Value = "  +7.16%  "
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "9.01"
$ws.Range("E12").Value = "  +14.15%  "
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").Value = "2.639.11"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "0.857"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "2.298.16"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "43.746.80"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").Value = "7.05"
$ws.Range("E20").Value = "  +13.37%  "
$ws.Range("D21").Value = "72.33"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "2.44"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "233.06"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "9.83"
$ws.Range("E24").Value = "  +7.07%  "
$ws.Range("D25").Value = "2.96"
$ws.Range("E25").Value = "  +7.97%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "11.64"
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("D28").Value = "43.24"
$ws.Range("E28").Value = "  +13.21%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "3.38"
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("D32").Value = "175.47"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").Value = "0.0936"
$ws.Range("E33").Value = "  +4.85%  "
$ws.Range("D34").Value = "21.64"
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("D35").Value = "5.71"
$ws.Range("E35").Value = "  +4.56%  "
$ws.Range("D36").Value = "4.82"
$ws.Range("E36").Value = "  +2.23%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  +2.99%  "
$ws.Range("D39").Value = "3.86"
$ws.Range("E39").Value = "  +9.11%  "
$ws.Range("E40").Value = "  -3.21%  "
$ws.Range("D41").Value = "14.30"
$ws.Range("E41").Value = "  +17.70%  "
$ws.Range("D42").Value = "75.49"
$ws.Range("E42").Value = "  +17.03%  "
$ws.Range("E43").Value = "  +2.59%  "
$ws.Range("E44").Value = "  +2.58%  "
$ws.Range("D45").Value = "6.36"
$ws.Range("E45").Value = "  +21.63%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "1.41"
$ws.Range("E47").Value = "  +2.44%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("E49").Value = "  +2.88%  "
$ws.Range("E50").Value = "  +3.52%  "
$ws.Range("D51").Value = "0.0992"
$ws.Range("E51").Value = "  -3.07%  "
